$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9972748160362244
$ws.Range("B1").Value = 2.066300630569458
$ws.Range("C1").Value = 9.404514312744141
$ws.Range("D1").Value = 2.685243844985962
$ws.Range("E1").Value = 1.403477072715759
